# Update "Förändrad" (Changed) date column (C) for rows 2-28 from 45509 (2024-08-05)
# to 45510 (2024-08-06), matching an automatic daily refresh of the log.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 28; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45509) {
        $cell.Value = 45510
    }
}
